$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contains a daily price log (rows 2..228) sorted with the newest
# entry always inserted at the top of the "Repollo" block (row 138). A new
# day of prices needs to be recorded, which pushes every existing row in the
# block (138..228) down by one, with the previously-last row (228) becoming
# the new last row (229). Insert a fresh row at 138 to reproduce that shift.
$ws.Rows("138").Insert()

# Populate the newly inserted row 138 with the new day's reading.
$ws.Cells.Item(138, 1).Value = 5
$ws.Cells.Item(138, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(138, 3).Value = "Maule"
$ws.Cells.Item(138, 4).Value = 44603
$ws.Cells.Item(138, 5).Value = 7
$ws.Cells.Item(138, 6).Value = 100112006
$ws.Cells.Item(138, 7).Value = "Repollo"
$ws.Cells.Item(138, 8).Value = "Crespo record"
$ws.Cells.Item(138, 9).Value = "Segunda"
$ws.Cells.Item(138, 10).Value = 2000
$ws.Cells.Item(138, 11).Value = 800
$ws.Cells.Item(138, 12).Value = 800
$ws.Cells.Item(138, 13).Value = 800
$ws.Cells.Item(138, 14).Value = "`$/unidad"
$ws.Cells.Item(138, 15).Value = "Región del Maule"
$ws.Cells.Item(138, 16).Value = 800
$ws.Cells.Item(138, 17).Value = 1
$ws.Cells.Item(138, 18).Value = "Hortaliza"
